$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "riqi" (date) column H for all data rows (2-41) to the new date.
# A leading apostrophe forces Excel to keep the value as text rather than
# reinterpreting the "yyyy-mm-dd" looking string as a real date serial.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 8).Value = "'2025-05-30"
}

# Fix the venue number in E41 from 2040 to 1040 to match the sequence,
# keeping it as text (same apostrophe trick so it doesn't become a number).
$ws.Cells.Item(41, 5).Value = "'1040"
